$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price & 1h volume change) per diff

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.932.36'
$ws.Range('E2').Value = '  -2.86%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.658.63'
$ws.Range('E3').Value = '  +0.46%  '

$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.71'
$ws.Range('E5').Value = '  -1.63%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.69'
$ws.Range('E6').Value = '  -7.70%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.598'
$ws.Range('E7').Value = '  -2.87%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('E9').Value = '  -3.74%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.59'
$ws.Range('E10').Value = '  -2.73%  '

$ws.Range('E11').Value = '  -1.54%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.10'
$ws.Range('E12').Value = '  -4.63%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.080.51'
$ws.Range('E13').Value = '  +1.07%  '

$ws.Range('E14').Value = '  +0.59%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.670.22'
$ws.Range('E15').Value = '  -0.03%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.933'
$ws.Range('E16').Value = '  -1.81%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.10'
$ws.Range('E17').Value = '  -2.48%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.968.65'
$ws.Range('E18').Value = '  -3.85%  '

$ws.Range('E19').Value = '  -2.49%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.82'
$ws.Range('E20').Value = '  -1.09%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.87'
$ws.Range('E21').Value = '  -5.68%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.34'
$ws.Range('E22').Value = '  +1.90%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '281.99'
$ws.Range('E23').Value = '  +2.35%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.06'
$ws.Range('E24').Value = '  -2.84%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '31.13'
$ws.Range('E25').Value = '  +0.77%  '

$ws.Range('E26').Value = '  +0.63%  '

$ws.Range('E27').Value = '  -0.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.58'
$ws.Range('E28').Value = '  -2.58%  '

$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '38.55'
$ws.Range('E29').Value = '  -7.07%  '

$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.17'
$ws.Range('E30').Value = '  -6.13%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.27'
$ws.Range('E31').Value = '  -1.54%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.76'
$ws.Range('E32').Value = '  -1.46%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.35'
$ws.Range('E33').Value = '  +1.87%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '156.01'
$ws.Range('E34').Value = '  +2.02%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0841'
$ws.Range('E35').Value = '  -2.57%  '

$ws.Range('E36').Value = '  -2.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.122'
$ws.Range('E37').Value = '  -2.52%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.06'
$ws.Range('E38').Value = '  +11.15%  '

$ws.Range('E39').Value = '  -0.90%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.06'
$ws.Range('E40').Value = '  -6.88%  '

$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').Value = '  -3.76%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0327'
$ws.Range('E42').Value = '  -2.15%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.98'
$ws.Range('E43').Value = '  -8.88%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.161.99'
$ws.Range('E44').Value = '  +3.59%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '94.17'
$ws.Range('E46').Value = '  -5.31%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '111.93'
$ws.Range('E47').Value = '  -3.24%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.32'
$ws.Range('E48').Value = '  -0.07%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.928.81'
$ws.Range('E49').Value = '  +1.01%  '

$ws.Range('E50').Value = '  -2.83%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.73'
$ws.Range('E51').Value = '  -8.00%  '
